# Workbook / worksheet handles
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the two derived columns (F: 产销率, G: 销售量) that are no longer tracked.
$ws.Range("F1:G81").EntireColumn.Delete()

# 2) The B/C quarter rows within every "<year>年" block of 4 rows (A,B,C,D) were
#    swapped in place (labels travel with their own data). Data rows run from
#    row 2 (2000年A) through row 81 (2019年D) in blocks of 4.
for ($blockStart = 2; $blockStart -le 81; $blockStart += 4) {
    $rowB = $blockStart + 1
    $rowC = $blockStart + 2

    for ($col = 1; $col -le 5; $col++) {
        $cellB = $ws.Cells.Item($rowB, $col)
        $cellC = $ws.Cells.Item($rowC, $col)

        $valB = $cellB.Value()
        $valC = $cellC.Value()

        $cellB.Value = $valC
        $cellC.Value = $valB
    }
}
